# Data Driver Framework Update
# Adds a new "Register" worksheet (after the existing "Login" sheet) that
# holds a small first-name / last-name / phone data table, and makes it the
# active sheet/tab (mirrors the workbook.xml activeTab + sheetView
# tabSelected flips in the target diff).

$wb = $excel.ActiveWorkbook

# Insert the new sheet right after "Login" so it becomes sheet #2.
$loginSheet = $wb.Worksheets.Item("Login")
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $loginSheet)
$ws2.Name = "Register"
$ws2.Activate()

# Row 1: fname1 / lname1 / PH1
$ws2.Range("A1").Value = "fname1"
$ws2.Range("B1").Value = "lname1"
$ws2.Range("C1").Value = "PH1"

# Row 2: fname2 / lname2 / PH2
$ws2.Range("A2").Value = "fname2"
$ws2.Range("B2").Value = "lname2"
$ws2.Range("C2").Value = "PH2"

# Match the sheet's stored selection (B3) like the other sheet in the book.
$ws2.Range("B3").Select()
